$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $result = $d.Content.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $result) {
        Write-Host "WARNING: replace failed for: $old"
    }
}

# Title
Replace-Text "Unveiling the Enigmatic Realm of Dark Matter" "Exploring the Realm of Science: Unraveling the Mysteries of the Natural World"

# Author name
Replace-Text " Amelia White" " Jonathan Smith"

# Email address: "amelia" + "." + "white@celestialobservatory" + "." + "org"
# becomes: "drjonathansmith@realscience" + "." + "org" (middle two runs removed)
Replace-Text "amelia.white@celestialobservatory" "drjonathansmith@realscience"

# Body paragraph 1 (first content paragraph)
Replace-Text "For centuries, humans have marveled at the vast expanse of the cosmos, seeking answers to the mysteries that lie beyond our visible world" "From the dawn of civilization, humans have sought to understand the complexities of the natural world"
Replace-Text " One of the most perplexing enigmas in modern cosmology is the existence of dark matter, an invisible substance that exerts a gravitational influence on galaxies and cosmic structures" " Science, a systematic and methodical approach to knowledge acquisition, has emerged as a powerful tool to decipher the mysteries that surround us"
Replace-Text " Its elusive nature has captivated scientists and fueled a quest to understand its properties, origin, and implications for our comprehension of the universe" " In its pursuit of truth, science unravels the interconnectedness of life, matter, and energy, demystifying phenomena and illuminating our place within the vast tapestry of existence"
Replace-Text "In the depths of space, galaxies spin with a peculiar velocity, suggesting the presence of unseen mass" "Embarking on this scientific journey, we delve into the realm of mathematics, where patterns, structures, and relationships unveil hidden truths"
Replace-Text " Gravitational lensing observations reveal the distortion of light around galaxies and clusters, hinting at the existence of a substantial amount of matter beyond what is visible" " Through formulas and equations, we unravel the intricate web of numerical relationships, unlocking the secrets of quantity, shape, and change"
Replace-Text " Additionally, the cosmic microwave background, a remnant radiation from the early universe, exhibits temperature fluctuations that are best explained by the influence of dark matter" " The language of mathematics empowers us to quantify, analyze, and predict natural phenomena, bridging the gap between abstract symbols and tangible realities"
Replace-Text "The search for dark matter particles has intensified in recent decades" "Venturing into the realm of chemistry, we explore the interactions between substances, delving into the intricacies of atomic structures, molecular bonds, and chemical reactions"
Replace-Text " Underground laboratories, shielded from cosmic rays and other interfering signals, have been constructed to detect the faint interactions of dark matter with ordinary matter" " We uncover the principles governing the composition, properties, and behavior of matter, revealing the fundamental building blocks of the universe"

# This replace also absorbs the deleted ". While these efforts...dark matter" run pair,
# keeping the final standalone "." run intact.
Replace-Text " Experiments utilizing sensitive detectors aim to uncover the elusive particles that may constitute dark matter, such as weakly interacting massive particles (WIMPs) or axions. While these efforts have yet to yield definitive results, they continue to push the boundaries of our knowledge and bring us closer to unraveling the secrets of dark matter" " Chemistry enables us to manipulate substances, synthesize new materials, and unravel the mysteries of life itself"

# Summary heading text unchanged ("Summary"), only font changes (handled globally below)

# Summary paragraph
Replace-Text "Dark matter, an enigmatic substance that pervades the universe, exerts a gravitational influence on galaxies and cosmic structures" "The exploration of science opens doors to a realm of wonder, where the mysteries of the natural world are gradually unveiled"
Replace-Text " Its existence is inferred through various observations, including the peculiar rotational velocities of galaxies, gravitational lensing effects, and the temperature fluctuations in the cosmic microwave background" " Through mathematics, we decipher patterns and relationships, quantifying and analyzing phenomena"
Replace-Text " The search for dark matter particles is an ongoing endeavor, with underground laboratories and sophisticated experiments attempting to detect these elusive particles" " Chemistry delves into the interactions between substances, revealing the fundamental building blocks of matter and unlocking the secrets of chemical reactions"
Replace-Text " The unraveling of the dark matter mystery promises to revolutionize our understanding of the universe and its composition" " These disciplines empower us to understand the complexities of the universe, inspiring us to seek knowledge and make meaningful contributions to the advancement of human understanding"

# Update font for the whole document: TimesNewToman -> Times New Roman
$full = $d.Range(0, $d.Content.End)
$full.Font.Name = "Times New Roman"

# Add a trailing empty paragraph at the end of the document body
$d.Paragraphs.Add() | Out-Null

Write-Host "Edit complete"
